# [RowConcealer] add new service "RowConcealer" to hide rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row of data just below the existing data ("this row remains" in A4)
$ws.Range("A5").Value = "this row will be hidden"

# Move the active selection down to the next empty row, as Excel does after
# typing a value into a cell.
$ws.Range("A6").Select()
